$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-19 Thursday" "2024-09-20 Friday"

Replace-Text "90×89=8010" "75×18=1350"
Replace-Text "81×72=5832" "54×57=3078"
Replace-Text "73×51=3723" "44×78=3432"
Replace-Text "14×69=966" "73×17=1241"
Replace-Text "64×80=5120" "74×94=6956"

Replace-Text "77×22=1694" "48×39=1872"
Replace-Text "33×73=2409" "79×49=3871"
Replace-Text "18×59=1062" "83×76=6308"
Replace-Text "29×82=2378" "34×19=646"
Replace-Text "69×76=5244" "26×76=1976"

Replace-Text "97×80=7760" "73×21=1533"
Replace-Text "49×42=2058" "47×73=3431"
Replace-Text "93×12=1116" "60×29=1740"
Replace-Text "35×52=1820" "65×59=3835"
Replace-Text "18×51=918" "66×21=1386"

Replace-Text "67×16=1072" "77×27=2079"
Replace-Text "95×78=7410" "27×69=1863"
Replace-Text "16×84=1344" "18×84=1512"
Replace-Text "92×14=1288" "74×63=4662"
Replace-Text "22×63=1386" "86×19=1634"

Replace-Text "61×73=4453" "50×80=4000"
Replace-Text "67×85=5695" "16×83=1328"
Replace-Text "74×36=2664" "33×62=2046"
Replace-Text "64×36=2304" "21×47=987"
Replace-Text "80×77=6160" "60×76=4560"
